$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Zero out the "target_central" block (rows 15-29)
$ws.Range("C15:C29").Value = 0

# Zero out the "target_linear" block (rows 43-57)
$ws.Range("C43:C57").Value = 0

# Zero out the "target_ambitious" block (rows 71-85)
$ws.Range("C71:C85").Value = 0

# Update selection to match the final state seen in the saved file
$ws.Activate()
$ws.Range("K37").Select()
